# "Add files via upload" - updates to curso_elton.xlsx:
#  - Plan1!F31 comment text rewritten (paragraphs merged, blank lines removed)
#  - Three new review comments added on Plan1!F32, F33, F34
#  - New rows of data (sim/nao + message) filled in on Plan1 C32:G34,
#    pulling in two new quote strings
#  - Sheet view refreshed (zoom + active selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# --- 1. Rewrite the F31 comment: the two trailing paragraphs get folded
#        into the second run and the blank separator lines between
#        paragraphs are removed. ---
$f31Text = "Alan Jose Nascimento:`n" + `
    "`n" + `
    "A pessoa mais indicada para dizer o que ele disse era ele.`n" + `
    "A pessoa mais indicada para dizer o que voce tem para dizer é você`n" + `
    "O cliente nao tem condiçoes de tomar a melhor decisão sem a minha ajuda`n" + `
    "Temos que conduzir o cliente pelo processo de decisão fazendo o papel de ocnselheiro.`n" + `
    "O meu papel e ajudar as pessoas a decidir. A melhor opção para as pessoas sou eu!!!!!`n" + `
    "Ao asusmir o papel de conselheiro temos mais coragem de mostrar as coisas para o cliente.`n"

$ws.Range("F31").Comment.Text($f31Text)

# --- 2. Add the three new comments (same reviewer, new notes) ---
$f32Text = "Alan Jose do Nascimento:`n`n"
$ws.Range("F32").AddComment($f32Text)

$f33Text = "Alan Jose do Nascimento:`n" + `
    "`n" + `
    "Pense em caixinhas e de preferencia sequencia:`n" + `
    "3 coisas que impedem as vendas`n" + `
    "2 amigos estavam na floretas e aparecue um urso e um dos amigos passou a amarrar os sapatos....`n" + `
    "3 coisas pode matar o resultado`n" + `
    "pense de forma estrategica`n" + `
    "1 - não oferecer, ou não oferecer o suficiente;`n" + `
    "2 - ofertar/oferecer do jeito errado`n" + `
    "3 - Estar oferecendo a coisa errada`n" + `
    "A solucao é a forma com o que você faz o seu cliente alcance o objetivo que ele tem.`n" + `
    "(O cliente só quer o beneficio)`n" + `
    "O cliente precisa comprar o beneficio, transformacao, o futuro que vc ta prometendo.`n" + `
    "alinhar volume de vendas x volume de pessoas suficiente`n" + `
    "ofertar o que as pessoas querem.!`n" + `
    "Pessoas diferentes querem coisas diferentes, por isso a escolha do nicho é importante.`n" + `
    "Prmissas:`n" + `
    "foca no publico;`n" + `
    "faça oferta no nivel suficiente;`n" + `
    "faça oferta do jeito certo;`n" + `
    "oferte para pessoas aquilo que elas querem (beneficio, transformacao, um futuro melhor)`n" + `
    "resolva as coisas por etapa`n"
$ws.Range("F33").AddComment($f33Text)

$f34Text = "Alan Jose do Nascimento:`n" + `
    "`n" + `
    "venda é diferente de oferta`n" + `
    "vc precisa estar se vendendo o tempo todo.`n" + `
    "Vendedor sempre vende e faz oferta de vez em quando.`n"
$ws.Range("F34").AddComment($f34Text)

# --- 3. Fill in the new data rows (32-34): "sim"/"nao" flags and the
#        message text in column G. Rows 32 and 33 re-use existing quotes;
#        33 and 34 introduce two brand-new quote strings. ---
$ws.Range("C32").Value = "sim"
$ws.Range("D32").Value = "não"
$ws.Range("G32").Value = "Você é a pessoa mais para opinar sobre a decisão de compra do seu cliente."

$ws.Range("C33").Value = "sim"
$ws.Range("D33").Value = "não"
$ws.Range("G33").Value = "Os nossos resultados são proporcionais a nossa capacidade de nos comunicarmos com os outros e com nós mesmos. Anthony Robbins"

$ws.Range("C34").Value = "sim"
$ws.Range("D34").Value = "não"
$ws.Range("G34").Value = "O melhor momento para fazer uma oferta ou um pedido de casamento é quando a outra parte já esta esperando, de preferencia desejando"

# --- 4. Refresh the sheet view: scrolled/zoomed in on the new rows, with
#        D35 as the active selection. ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 70
$ws.Range("D35").Select()
